# Update the monitoring leaderboard (data_monitoreo_la_peñita) with the
# latest "total_registros" counts. The sheet is sorted descending by
# total_registros, so a few names shift rows as their counts grow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (name, total_registros) pairs, in final sorted order for rows 2..12.
$rows = @(
    @{ Row = 2;  Name = "PEREZ VEGA ANA YSABEL";            Total = 142 },
    @{ Row = 3;  Name = "ZAPATA ZETA ROSA ARACELI";         Total = 137 },
    @{ Row = 4;  Name = "GARAVITO LEON IVONNE LISSETH";     Total = 134 },
    @{ Row = 5;  Name = "TIMOTEO BAYONA SHARYN LISSETH";    Total = 133 },
    @{ Row = 6;  Name = "PANTA MONZON SHIRLEY MARIBEL";     Total = 126 },
    @{ Row = 7;  Name = "NIÑO GUERRERO ANYELA MELINA";      Total = 111 },
    @{ Row = 8;  Name = "VALLE SILVA SUTMMER ORFELINDA";    Total = 108 },
    @{ Row = 9;  Name = "TIZON NUÑEZ FRESIA YAMILI";        Total = 100 },
    @{ Row = 10; Name = "CASTRO JUAREZ MARIA ISABEL";       Total = 99  },
    @{ Row = 12; Name = "CHERO JUAREZ ANYELA TATIANA";      Total = 79  }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    $ws.Cells.Item($r.Row, 2).Value = $r.Total
}
